$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.8647090984694306, 0.6279050292696604, 1.6680471595701123, 2.9181192133307254)"
$ws.Range("C2").Value = "JSU(-1.1927153770673833, 1.214860700748328, 2.7406797929057496, 5.08103349729418)"
$ws.Range("D2").Value = "NIG(1.0271281951356217, 0.8037764716274598, 0.42598094621036353, 3.1210032669151033)"
$ws.Range("E2").Value = "NIG(1.2642059964601333, 0.9860262200137587, 3.71168223775103, 4.965462266034573)"
